$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C13 capacitor entry -> now covers C13, C16 and C22 with a new part number,
# quantity bumped from 1 to 3, and it is no longer marked "Nicht bestücken" (DNP cleared).
$ws.Range("A2").Value = "C13,C16,C22"
$ws.Range("B2").Value = "MKS4C042204C00KSSD"
$ws.Range("C2").Value = "SamacSys_Parts:MKS4C042204C00KSSD"
$ws.Range("D2").Value = 3
$ws.Range("E2").ClearContents()

# Row 4: connector designator list grows from J1-J5 to J1-J11, quantity 5 -> 11.
$ws.Range("A4").Value = "J1,J2,J3,J4,J5,J6,J7,J8,J9,J10,J11"
$ws.Range("D4").Value = 11
